$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Day 7" row of data (row 8) to the tracker
$ws.Range("A8").Value = "Day 7"

$ws.Range("B7").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B8").Value = Get-Date -Year 2025 -Month 6 -Day 1 -Hour 0 -Minute 0 -Second 0

$ws.Range("C8").Value = "Koko Eating Bananas"
$ws.Range("C8").Interior.Color = 65535

$ws.Range("D8").Value = "Search in Rotated Sorted Array"
$ws.Range("E8").Value = "Find Minimum in Rotated Sorted Array"
$ws.Range("F8").Value = "Binary Search, Arrays"
$ws.Range("G8").Value = "S"
$ws.Range("H8").Value = "YES"
